$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the "Monday 2016-12-21" entry block to make room
# for the new "2016-12-26" trading-plan entry (data row + blank spacer row).
$ws.Range("A3:G4").EntireRow.Insert()

# --- New entry for 2016-12-26 / 2016-12-27 -------------------------------
$ws.Range("A3").Value = 20161226
$ws.Range("B3").Value = "Monday"
$ws.Range("C3").Value = 20161227
$ws.Range("D3").Value = "Tuesday"
$ws.Range("F3").Value = "Very bad decision buying without knowing the weather could change (turn colder in early Jan) and the fund is not settled. NEVER do it AGAIN! Do not fight the trend"
$ws.Range("E3").Value = "Already under the water and possibly already erased all last week's gain… find a chagne to get out to stop much losses. Two days away from Thursday's report, which could be a small withdraw due to last week's warm weather. If the NG pircce don't pass 3.777, we may have a chance to hold and make profits, becauase that means the turning point of the up-trend "

# --- Wrap-text formatting sweep over the whole table ----------------------
# Column A (LogTime) becomes wrap-text only (style index 1).
$ws.Range("A1").WrapText = $true
$ws.Range("A3").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Range("A9").WrapText = $true
$ws.Range("A11").WrapText = $true

# Columns B, C, D gain wrap-text on top of their existing left alignment
# (style index 2 -> 3) across every used row.
$ws.Range("B1:D11").WrapText = $true

# Column F keeps the same treatment, but only on the rows that actually hold
# an F cell (F1/F5/F7 pre-existing, F2/F3 newly introduced need both the
# left alignment and the wrap text set explicitly since they start blank).
$ws.Range("F1").WrapText = $true
$ws.Range("F2").HorizontalAlignment = -4131
$ws.Range("F2").WrapText = $true
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").WrapText = $true
$ws.Range("F5").WrapText = $true
$ws.Range("F7").WrapText = $true

# G2 / G3 pick up the same left-aligned (non-wrap) style already used by G1.
$ws.Range("G2").HorizontalAlignment = -4131
$ws.Range("G3").HorizontalAlignment = -4131

# --- Row heights (auto-fit result once wrap text covers the wider table) --
$ws.Rows(3).RowHeight = 57.6
$ws.Rows(5).RowHeight = 100.8
$ws.Rows(7).RowHeight = 43.2
$ws.Rows(9).RowHeight = 43.2
$ws.Rows(11).RowHeight = 43.2
